$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the incorrect "isSourceOf" / "isDerivedFrom" header labels from H1:I1,
# leaving the cells empty but keeping their existing (header) formatting.
$ws.Range("H1:I1").ClearContents()

# Reflect the resulting selection (H1 was the active/selected cell after the edit).
$ws.Range("H1").Select()
